$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 41 (inherits number formats/styles from the row above,
# matching the existing "Baseline 2010-18 CNNN" summary rows).
$ws.Rows.Item(41).Insert()

$ws.Range("A41").Value2 = "CW3M"
$ws.Range("B41").Value2 = "Baseline 2010-18 C128"
$ws.Range("C41").Value2 = "2010-18"
$ws.Range("D41").Value2 = 1186.9773491111109
$ws.Range("E41").Value2 = 1901.5157334444443
$ws.Range("F41").Value2 = 0.97970299999999988
$ws.Range("G41").Value2 = 280.33542888888883
$ws.Range("H41").Value2 = 9.775355222222224
$ws.Range("I41").Value2 = 5.3870271111111121
$ws.Range("J41").Value2 = 8.145128999999999
$ws.Range("K41").Value2 = 645.93808322222219
$ws.Range("L41").Value2 = 83.47062044444445
$ws.Range("M41").Value2 = 1455.5792641111111
$ws.Range("N41").Value2 = 1191.1918266666667
$ws.Range("O41").Value2 = 4661.9885253333332
$ws.Range("P41").Value2 = 27227.338324888889
$ws.Range("Q41").Value2 = -0.64567288888888896
$ws.Range("R41").Value2 = -0.00020755555555555555
$ws.Range("S41").Value2 = "2010-18"

$ws.Range("U49").Select()
